# Bug 129: Image formatting does not accept non decimal numbers: fixed regex pattern
#
# The single sample paragraph is rewritten into four paragraphs that
# exercise integer / decimal / incomplete-decimal image size parameters.

$d = $word.ActiveDocument

# Helper: insert each string in $parts as its own separate run (all sharing the
# paragraph's current/default run formatting) at the *start* of paragraph
# number $paraIndex. We do this by typing each part, then temporarily
# splitting the paragraph with a paragraph mark (so the just-typed text
# becomes - and stays - its own run), and finally deleting those temporary
# paragraph marks again to merge everything back into one paragraph while
# keeping the run boundaries that were created.
function Insert-Runs($paraIndex, $parts) {
    $insPoint = $d.Paragraphs($paraIndex).Range
    $insPoint.Collapse(1)
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $insPoint.InsertAfter($parts[$i])
        $insPoint.MoveStart(1, $parts[$i].Length)
        if ($i -lt $parts.Length - 1) {
            $insPoint.InsertParagraphAfter()
            $insPoint.MoveStart(1, 1)
        }
    }
    $mergeCount = $parts.Length - 1
    for ($i = 0; $i -lt $mergeCount; $i++) {
        $p = $d.Paragraphs($paraIndex)
        $endOfPara = $p.Range.End
        $delRange = $d.Range($endOfPara - 1, $endOfPara)
        $delRange.Delete()
    }
}

# Start from a single empty paragraph (keeps the original paragraph mark's
# properties - language etc. - intact).
$d.Content.Delete()

# Paragraph 1: "Size parameters as integers {image:w4cm;h4cm }"
Insert-Runs 1 @("Size parameters as integers", " ", "{image:w4cm;h4cm }")

# Paragraph 2: "Size parameters as decimals {image:w3.5cm;h3.5cm}"
$d.Paragraphs(1).Range.InsertParagraphAfter()
Insert-Runs 2 @("Size parameters as decimals", " ", "{image", ":w", "3", ".5cm;h", "3", ".5cm", "}")

# Paragraph 3: "Size parameters as incomplete decimals {image:w3.cm;h3.cm}"
$d.Paragraphs(2).Range.InsertParagraphAfter()
Insert-Runs 3 @("Size parameters as ", "incomplete ", "decimals {image:w3.cm;h3.cm}")

# Paragraph 4: trailing empty paragraph.
$d.Paragraphs(3).Range.InsertParagraphAfter()

Write-Output ("Paragraph count: " + $d.Paragraphs.Count)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Output ("  " + $i + ": [" + $d.Paragraphs($i).Range.Text + "]")
}
